$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "65.532.20"
Set-TextValue "E2" "  -1.04%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.298.29"
Set-TextValue "E3" "  -0.16%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.29%  "

# Row 5 - BNB
Set-TextValue "D5" "579.84"
Set-TextValue "E5" "  +4.16%  "

# Row 6 - Solana
Set-TextValue "D6" "182.99"
Set-TextValue "E6" "  -3.49%  "

# Row 7 - USDC
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.06%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.293.89"
Set-TextValue "E8" "  -0.08%  "

# Row 9 - XRP
Set-TextValue "D9" "0.571"
Set-TextValue "E9" "  -2.62%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.176"
Set-TextValue "E10" "  -4.61%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.572"
Set-TextValue "E11" "  -2.76%  "

# Row 12 - Avalanche
Set-TextValue "D12" "46.58"
Set-TextValue "E12" "  -2.32%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000264"
Set-TextValue "E13" "  -2.65%  "

# Row 14 - BitcoinCash
Set-TextValue "D14" "635.61"
Set-TextValue "E14" "  +3.20%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.823.53"
Set-TextValue "E15" "  -0.21%  "

# Row 16 - Polkadot
Set-TextValue "D16" "8.43"
Set-TextValue "E16" "  -2.70%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "65.651.92"
Set-TextValue "E17" "  -0.79%  "

# Row 18 - TRON
Set-TextValue "E18" "  +0.11%  "

# Row 19 - Chainlink
Set-TextValue "D19" "17.70"
Set-TextValue "E19" "  -1.99%  "

# Row 20 - WrappedEther
Set-TextValue "D20" "3.295.21"
Set-TextValue "E20" "  -0.01%  "

# Row 21 - Uniswap
Set-TextValue "D21" "10.97"
Set-TextValue "E21" "  -0.66%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.888"
Set-TextValue "E22" "  -2.29%  "

# Row 23 - InternetComputer(DFINITY)
Set-TextValue "D23" "17.91"
Set-TextValue "E23" "  -2.68%  "

# Row 24 - Litecoin
Set-TextValue "D24" "100.75"
Set-TextValue "E24" "  -1.83%  "

# Row 25 - Toncoin
Set-TextValue "D25" "4.95"
Set-TextValue "E25" "  -0.12%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "3.96"
Set-TextValue "E26" "  +0.58%  "

# Row 27 - ImmutableX
Set-TextValue "D27" "2.74"
Set-TextValue "E27" "  -0.28%  "

# Row 28 - RenderToken
Set-TextValue "D28" "9.37"
Set-TextValue "E28" "  -2.88%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "30.82"
Set-TextValue "E29" "  +1.78%  "

# Row 30 - Filecoin
Set-TextValue "D30" "8.36"
Set-TextValue "E30" "  -3.77%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "6.51"
Set-TextValue "E31" "  -1.18%  "

# Row 32 - Bittensor
Set-TextValue "D32" "580.11"
Set-TextValue "E32" "  +2.63%  "

# Row 33 - dogwifhat
Set-TextValue "E33" "  -10.48%  "

# Row 34 - Cosmos
Set-TextValue "D34" "10.87"
Set-TextValue "E34" "  -1.91%  "

# Row 35 - Maker
Set-TextValue "D35" "3.842.03"
Set-TextValue "E35" "  +1.00%  "

# Row 36 - Hedera
Set-TextValue "D36" "0.105"
Set-TextValue "E36" "  -0.60%  "

# Row 37 - Dai (unchanged)

# Row 38 - OKB
Set-TextValue "D38" "55.65"
Set-TextValue "E38" "  -3.18%  "

# Row 39 - Kaspa
Set-TextValue "E39" "  -3.03%  "

# Row 40 - ApeXProtocol
Set-TextValue "D40" "3.41"
Set-TextValue "E40" "  +5.22%  "

# Row 41 - InjectiveProtocol
Set-TextValue "D41" "32.42"
Set-TextValue "E41" "  -5.33%  "

# Row 42 - PEPE
Set-TextValue "D42" "0.0₃0688"
Set-TextValue "E42" "  -5.18%  "

# Row 43 - Stacks
Set-TextValue "E43" "  -6.23%  "

# Row 44 - Fetch.AI
Set-TextValue "E44" "  -5.53%  "

# Row 45 - TheGraph
Set-TextValue "D45" "0.332"
Set-TextValue "E45" "  -1.84%  "

# Row 46 - VeChain
Set-TextValue "D46" "0.0405"
Set-TextValue "E46" "  -3.87%  "

# Row 47 - CoreDAO
Set-TextValue "D47" "3.05"
Set-TextValue "E47" "  -4.86%  "

# Row 48 - now Stellar (was FirstDigitalUSD)
Set-TextValue "B48" "Stellar"
Set-TextValue "C48" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D48" "0.127"
Set-TextValue "E48" "  -2.05%  "

# Row 49 - now FirstDigitalUSD (was Stellar)
Set-TextValue "B49" "FirstDigitalUSD"
Set-TextValue "C49" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D49" "1.00"
Set-TextValue "E49" "  +0.42%  "

# Row 50 - ThetaToken
Set-TextValue "E50" "  -2.39%  "

# Row 51 - Monero
Set-TextValue "D51" "129.78"
Set-TextValue "E51" "  +5.70%  "
